# Update status of "Search Vehicles" (D7) and "Select Vehicle" (D8)
# from "InProgress" to "Completed" on the "Inc Matrix" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inc Matrix")

$ws.Range("D7").Value = "Completed"
$ws.Range("D8").Value = "Completed"

# Recolor the cells to match the "Completed" status fill (green),
# matching the look of the other "Completed" cells (e.g. D2:D6).
$ws.Range("D7:D8").Interior.Color = $ws.Range("D2").Interior.Color

# Move the active selection to D14, matching the saved cursor position.
$ws.Range("D14").Select()
